$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Params")

# Duplicate row 5 into the new row 6 (same values/number formats/styles),
# then change column A to the new subject "SSTest"
$ws.Range("A5:AW5").Copy($ws.Range("A6:AW6")) | Out-Null

$ws.Range("A6").Value = "SSTest"

# The Copy() above turns #N/A error cells into plain numeric error codes;
# restore them as real #N/A errors to match row 5.
$errCols = @("D", "AL", "AM", "AN", "AO", "AP", "AQ", "AR", "AS", "AT", "AU", "AV", "AW")
foreach ($col in $errCols) {
    $ws.Range($col + "6").Value = "#N/A"
}

$ws.Range("A7").Select() | Out-Null
